# Update weekly price data rows (7-16) on the active sheet.
# Each row's values are shifted from the following week's figures
# (row 16 receives the newest week's data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("D7").Value = 44482
$ws.Range("J7").Value = 160
$ws.Range("K7").Value = 7000
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = 7500
$ws.Range("P7").Value = 375

# Row 8
$ws.Range("D8").Value = 44755

# Row 9
$ws.Range("D9").Value = 44643
$ws.Range("K9").Value = 8000
$ws.Range("L9").Value = 9000
$ws.Range("M9").Value = 8500
$ws.Range("P9").Value = 425

# Row 10
$ws.Range("D10").Value = 44358
$ws.Range("K10").Value = 7500
$ws.Range("L10").Value = 8000
$ws.Range("M10").Value = 7750
$ws.Range("P10").Value = 388

# Row 11
$ws.Range("D11").Value = 44435
$ws.Range("J11").Value = 302
$ws.Range("K11").Value = 7000
$ws.Range("M11").Value = 7500
$ws.Range("P11").Value = 375

# Row 12
$ws.Range("D12").Value = 44162
$ws.Range("J12").Value = 50
$ws.Range("K12").Value = 8000
$ws.Range("M12").Value = 8000
$ws.Range("P12").Value = 400

# Row 13
$ws.Range("D13").Value = 44295
$ws.Range("J13").Value = 70

# Row 14
$ws.Range("D14").Value = 44273

# Row 15
$ws.Range("D15").Value = 44650
$ws.Range("J15").Value = 160
$ws.Range("K15").Value = 9000
$ws.Range("L15").Value = 10000
$ws.Range("M15").Value = 9500
$ws.Range("P15").Value = 475

# Row 16
$ws.Range("D16").Value = 45035
$ws.Range("J16").Value = 70
$ws.Range("L16").Value = 9000
$ws.Range("M16").Value = 9000
$ws.Range("P16").Value = 450
